# Regenerate save_data column G ("K" = strikeouts) to use actual strikeout
# counts instead of the previous "Strike#" (number of pitches that were
# strikes) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values for rows 2-28, replacing the old Strike# values.
$kValues = @{
    2  = 0
    3  = 1
    4  = 6
    5  = 5
    6  = 7
    7  = 5
    8  = 5
    9  = 7
    10 = 4
    11 = 12
    12 = 10
    13 = 6
    14 = 8
    15 = 6
    16 = 4
    17 = 8
    18 = 6
    19 = 2
    20 = 5
    21 = 3
    22 = 5
    23 = 4
    24 = 7
    25 = 4
    26 = 5
    27 = 3
    28 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
